$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-18 Saturday" "2025-01-19 Sunday"

Replace-Text "228÷4=57, 0" "134÷9=14, 8"
Replace-Text "165÷4=41, 1" "762÷4=190, 2"
Replace-Text "611÷3=203, 2" "339÷8=42, 3"
Replace-Text "245÷4=61, 1" "607÷6=101, 1"
Replace-Text "154÷5=30, 4" "799÷8=99, 7"

Replace-Text "531÷6=88, 3" "491÷7=70, 1"
Replace-Text "301÷8=37, 5" "884÷7=126, 2"
Replace-Text "891÷2=445, 1" "918÷7=131, 1"
Replace-Text "163÷7=23, 2" "779÷8=97, 3"
Replace-Text "893÷6=148, 5" "813÷3=271, 0"

Replace-Text "436÷2=218, 0" "799÷9=88, 7"
Replace-Text "498÷8=62, 2" "713÷3=237, 2"
Replace-Text "415÷7=59, 2" "360÷9=40, 0"
Replace-Text "216÷4=54, 0" "837÷6=139, 3"
Replace-Text "119÷5=23, 4" "886÷7=126, 4"

Replace-Text "794÷4=198, 2" "584÷5=116, 4"
Replace-Text "627÷6=104, 3" "638÷9=70, 8"
Replace-Text "327÷4=81, 3" "681÷9=75, 6"
Replace-Text "337÷7=48, 1" "806÷8=100, 6"
Replace-Text "780÷7=111, 3" "578÷9=64, 2"

Replace-Text "716÷2=358, 0" "824÷5=164, 4"
Replace-Text "630÷7=90, 0" "844÷9=93, 7"
Replace-Text "562÷9=62, 4" "775÷9=86, 1"
Replace-Text "301÷7=43, 0" "913÷2=456, 1"
Replace-Text "425÷4=106, 1" "661÷7=94, 3"
